# Added logic to update duplicated values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 2 and 3), pushing the
# existing Layout_* rows down.
$ws.Range("A2:B3").Insert()

$ws.Range("A2").Value = "Year"
$ws.Range("B2").Value = 2015

$ws.Range("A3").Value = "Term"
$ws.Range("B3").Value = "Spring"

$ws.Range("B3").Select()
